# Update automatico via Actualizar 04-16-2021 12-45-58
#
# The sheet keeps a rolling history of availability checks: each "batch"
# of 14 service rows shares one timestamp in column D. On every run the
# whole history shifts forward by one slot - the newest batch (rows 2-15)
# gets the freshly captured timestamp, and the two older batches (rows
# 16-29 and 30-43) are overwritten with the timestamp that used to belong
# to the batch above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newestTimestamp = 44302.53166146284
$middleTimestamp = 44302.51026142361
$oldestTimestamp = 44302.48887386574

for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 4).Value = $newestTimestamp
}

for ($row = 16; $row -le 29; $row++) {
    $ws.Cells.Item($row, 4).Value = $middleTimestamp
}

for ($row = 30; $row -le 43; $row++) {
    $ws.Cells.Item($row, 4).Value = $oldestTimestamp
}
